$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 24.71000000000042
$ws.Range("H2").Value = 0.0000000006376315031531021
$ws.Range("I2").Value = 0.0000000006376315031531021
$ws.Range("L2").Value = 41.05457925348692
$ws.Range("M2").Value = "[29.58780014514555, 52.52135836182829]"
$ws.Range("N2").Value = 0.000000004960717170376938
$ws.Range("O2").Value = 0.000000004960717170376938
$ws.Range("P2").Value = 1.440289725069195
$ws.Range("Q2").Value = "[1.1132370364071935, 1.7673424137311962]"
$ws.Range("R2").Value = 0.00000000001936739657537601
$ws.Range("S2").Value = 0.00000000001936739657537601
$ws.Range("T2").Value = 48.65416199422457
$ws.Range("U2").Value = "[41.59920742881233, 55.709116559636804]"
$ws.Range("X2").Value = 19.04574574574607
$ws.Range("Y2").Value = 17.75953953953984
$ws.Range("Z2").Value = 20.3319519519523

# Row 3
$ws.Range("F3").Value = 24.71000000000042
$ws.Range("H3").Value = 0.00000001131946814059148
$ws.Range("I3").Value = 0.00000001131946814059148
$ws.Range("L3").Value = 43.25520176751031
$ws.Range("M3").Value = "[28.804433746660855, 57.70596978835976]"
$ws.Range("N3").Value = 0.0000002829847882690473
$ws.Range("O3").Value = 0.0000002829847882690473
$ws.Range("P3").Value = 1.33965812855781
$ws.Range("Q3").Value = "[0.9874475407679633, 1.6918687163476562]"
$ws.Range("R3").Value = 0.000000001079211608967512
$ws.Range("S3").Value = 0.000000001079211608967512
$ws.Range("T3").Value = 58.53601991229636
$ws.Range("U3").Value = "[50.24768930264294, 66.82435052194978]"
$ws.Range("X3").Value = 19.44150150150183
$ws.Range("Y3").Value = 18.05635635635667
$ws.Range("Z3").Value = 20.826646646647

# Row 4
$ws.Range("F4").Value = 24.71000000000042
$ws.Range("H4").Value = 0.0000003577988073599769
$ws.Range("I4").Value = 0.0000003577988073599769
$ws.Range("L4").Value = 44.22132903320777
$ws.Range("M4").Value = "[25.659836720101517, 62.782821346314016]"
$ws.Range("N4").Value = 0.00001794072353833442
$ws.Range("O4").Value = 0.00001794072353833442
$ws.Range("P4").Value = 0.9748685912040393
$ws.Range("Q4").Value = "[0.5471843060306538, 1.4025528763774249]"
$ws.Range("R4").Value = 0.00003541999812650687
$ws.Range("S4").Value = 0.00003541999812650687
$ws.Range("T4").Value = 60.80733191472236
$ws.Range("U4").Value = "[51.087584140051476, 70.52707968939325]"
$ws.Range("V4").Value = 0.0000000000000002220446049250313
$ws.Range("W4").Value = 0.0000000000000002220446049250313
$ws.Range("X4").Value = 20.87611611611648
$ws.Range("Y4").Value = 19.19415415415449
$ws.Range("Z4").Value = 22.55807807807847

# Row 5
$ws.Range("F5").Value = 24.71000000000042
$ws.Range("H5").Value = 0.000001214393403392577
$ws.Range("I5").Value = 0.000001214393403392577
$ws.Range("L5").Value = 36.97235654734548
$ws.Range("M5").Value = "[20.719613948776583, 53.225099145914385]"
$ws.Range("N5").Value = 0.00003649899046220817
$ws.Range("O5").Value = 0.00003649899046220817
$ws.Range("P5").Value = 0.4842895582110387
$ws.Range("Q5").Value = "[0.03144737390980712, 0.9371317425122703]"
$ws.Range("R5").Value = 0.03663507051623638
$ws.Range("S5").Value = 0.03663507051623638
$ws.Range("T5").Value = 52.74330178242745
$ws.Range("U5").Value = "[44.243781927356494, 61.24282163749841]"
$ws.Range("V5").Value = 0.0000000000000004440892098500626
$ws.Range("W5").Value = 0.0000000000000004440892098500626
$ws.Range("X5").Value = 22.80542542542582
$ws.Range("Y5").Value = 21.02452452452489
$ws.Range("Z5").Value = 24.58632632632675

# Row 6
$ws.Range("F6").Value = 24.71000000000042
$ws.Range("H6").Value = 0.0000008446692743024897
$ws.Range("I6").Value = 0.0000008446692743024897
$ws.Range("L6").Value = 36.76141366851108
$ws.Range("M6").Value = "[20.365559043504255, 53.15726829351791]"
$ws.Range("N6").Value = 0.00004521618264519134
$ws.Range("O6").Value = 0.00004521618264519134
$ws.Range("P6").Value = 0.8993948938205012
$ws.Range("Q6").Value = "[0.45913165908319353, 1.3396581285578089]"
$ws.Range("R6").Value = 0.0001630325751333928
$ws.Range("S6").Value = 0.0001630325751333928
$ws.Range("T6").Value = 51.40646577900528
$ws.Range("U6").Value = "[43.01865126163197, 59.79428029637859]"
$ws.Range("V6").Value = 0.0000000000000004440892098500626
$ws.Range("W6").Value = 0.0000000000000004440892098500626
$ws.Range("X6").Value = 21.1729329329333
$ws.Range("Y6").Value = 19.44150150150184
$ws.Range("Z6").Value = 22.90436436436475

# Row 7
$ws.Range("B7").Value = 1
$ws.Range("F7").Value = 24.71000000000042
$ws.Range("H7").Value = 0.000000002499635276542733
$ws.Range("I7").Value = 0.000000002499635276542733
$ws.Range("L7").Value = 50.65484244173497
$ws.Range("M7").Value = "[34.07369742034825, 67.23598746312169]"
$ws.Range("N7").Value = 0.00000018500369303176
$ws.Range("O7").Value = 0.00000018500369303176
$ws.Range("P7").Value = 0.5471843060306538
$ws.Range("Q7").Value = "[0.20755266780472947, 0.8868159442565782]"
$ws.Range("R7").Value = 0.002219207133215439
$ws.Range("S7").Value = 0.002219207133215439
$ws.Range("T7").Value = 54.56902879263047
$ws.Range("U7").Value = "[45.625420829549086, 63.51263675571185]"
$ws.Range("V7").Value = 0.0000000000000004440892098500626
$ws.Range("W7").Value = 0.0000000000000004440892098500626
$ws.Range("X7").Value = 22.55807807807847
$ws.Range("Y7").Value = 21.22240240240276
$ws.Range("Z7").Value = 23.89375375375417

# Row 8
$ws.Range("B8").Value = 0
$ws.Range("F8").Value = 22.53000000000008
$ws.Range("H8").Value = 0.00000001579550712449418
$ws.Range("I8").Value = 0.00000001579550712449418
$ws.Range("L8").Value = 54.40781025099125
$ws.Range("M8").Value = "[37.39156757900345, 71.42405292297906]"
$ws.Range("N8").Value = 0.00000006926236251203477
$ws.Range("O8").Value = 0.00000006926236251203477
$ws.Range("P8").Value = 0.1069210712933462
$ws.Range("Q8").Value = "[-0.25786846606042335, 0.47171060864711567]"
$ws.Range("R8").Value = 0.5579154612048591
$ws.Range("S8").Value = 0.5579154612048591
$ws.Range("T8").Value = 58.53697839810313
$ws.Range("U8").Value = "[47.813852500486306, 69.26010429571996]"
$ws.Range("V8").Value = 0.00000000000002464695114667848
$ws.Range("W8").Value = 0.00000000000002464695114667848
$ws.Range("X8").Value = 22.14660660660669
$ws.Range("Y8").Value = 20.83855855855864
$ws.Range("Z8").Value = 23.45465465465474

# Row 9
$ws.Range("F9").Value = 22.53000000000008
$ws.Range("H9").Value = 0.000000001244478409212491
$ws.Range("I9").Value = 0.000000001244478409212491
$ws.Range("L9").Value = 43.8197440948878
$ws.Range("M9").Value = "[29.09163142994033, 58.54785675983526]"
$ws.Range("N9").Value = 0.0000003203989733524537
$ws.Range("O9").Value = 0.0000003203989733524537
$ws.Range("P9").Value = 0.4968685077749626
$ws.Range("Q9").Value = "[0.15723686954903826, 0.836500146000887]"
$ws.Range("R9").Value = 0.005075208778289264
$ws.Range("S9").Value = 0.005075208778289264
$ws.Range("T9").Value = 48.81305783645989
$ws.Range("U9").Value = "[41.11558605249015, 56.51052962042963]"
$ws.Range("V9").Value = 0.0000000000000002220446049250313
$ws.Range("W9").Value = 0.0000000000000002220446049250313
$ws.Range("X9").Value = 20.74834834834842
$ws.Range("Y9").Value = 19.53051051051058
$ws.Range("Z9").Value = 21.96618618618627

# Row 10
$ws.Range("F10").Value = 22.53000000000008
$ws.Range("H10").Value = 0.00000001164152496624382
$ws.Range("I10").Value = 0.00000001164152496624382
$ws.Range("L10").Value = 46.01382308444226
$ws.Range("M10").Value = "[29.072909629023187, 62.95473653986133]"
$ws.Range("N10").Value = 0.000001892917982360132
$ws.Range("O10").Value = 0.000001892917982360132
$ws.Range("P10").Value = 0.9371317425122703
$ws.Range("Q10").Value = "[0.5597632555945777, 1.3145002294299628]"
$ws.Range("R10").Value = 0.000009147476090820561
$ws.Range("S10").Value = 0.000009147476090820561
$ws.Range("T10").Value = 58.50141176581961
$ws.Range("U10").Value = "[49.89845986282772, 67.1043636688115]"
$ws.Range("V10").Value = 0
$ws.Range("W10").Value = 0
$ws.Range("X10").Value = 19.16966966966974
$ws.Range("Y10").Value = 17.81651651651658
$ws.Range("Z10").Value = 20.5228228228229

